$d = $word.ActiveDocument

# Classic Word color constants (OLE COLORREF, 0x00BBGGRR packed as R | G<<8 | B<<16)
$wdColorAutomatic = -16777216
$wdColorRed    = 192        # RGB(C0,00,00) -> C00000
$wdColorGreen  = 9293992    # RGB(A8,D0,8D) -> A8D08D
$wdTextureNone = 0

function Set-CellShading($table, $rowIdx, $colIdx, $fillColor) {
    $cell = $table.Cell($rowIdx, $colIdx)
    $shd = $cell.Shading
    $shd.Texture = $wdTextureNone
    $shd.ForegroundPatternColor = $wdColorAutomatic
    $shd.BackgroundPatternColor = $fillColor
}

# --- Table 1 ---
$t1 = $d.Tables.Item(1)
# Row 2 (output 0) -> red
Set-CellShading $t1 2 1 $wdColorRed
Set-CellShading $t1 2 2 $wdColorRed
# Rows 3-5 (output 1) -> green
Set-CellShading $t1 3 1 $wdColorGreen
Set-CellShading $t1 3 2 $wdColorGreen
Set-CellShading $t1 4 1 $wdColorGreen
Set-CellShading $t1 4 2 $wdColorGreen
Set-CellShading $t1 5 1 $wdColorGreen
Set-CellShading $t1 5 2 $wdColorGreen

# --- Table 2 ---
$t2 = $d.Tables.Item(2)
# Rows 2-4 (output 1) -> green
Set-CellShading $t2 2 1 $wdColorGreen
Set-CellShading $t2 2 2 $wdColorGreen
Set-CellShading $t2 3 1 $wdColorGreen
Set-CellShading $t2 3 2 $wdColorGreen
Set-CellShading $t2 4 1 $wdColorGreen
Set-CellShading $t2 4 2 $wdColorGreen
# Row 5 (output 0) -> red
Set-CellShading $t2 5 1 $wdColorRed
Set-CellShading $t2 5 2 $wdColorRed
